$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 74-82 (data beyond the new end of series)
$ws.Rows("74:82").Delete()

# Update forecasted/recomputed values in column B (rows 4-73)
$ws.Range("B4").Value = -0.9716635481116656
$ws.Range("B5").Value = 0.2805049688864443
$ws.Range("B6").Value = -0.6760650065702858
$ws.Range("B7").Value = -0.04573875985788579
$ws.Range("B8").Value = -0.4726963440578849
$ws.Range("B9").Value = 0.1769623127656392
$ws.Range("B10").Value = 0.2080098490516883
$ws.Range("B11").Value = 0.2409726375546795
$ws.Range("B12").Value = 0.03734130806860806
$ws.Range("B13").Value = 0.01272818083837047
$ws.Range("B14").Value = 0.8
$ws.Range("B15").Value = 0.3
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = 0.5
$ws.Range("B18").Value = 0.1
$ws.Range("B19").Value = 0.4
$ws.Range("B20").Value = 0.5
$ws.Range("B21").Value = 0.3
$ws.Range("B22").Value = 0.1
$ws.Range("B23").Value = 0.2
$ws.Range("B24").Value = 0.9
$ws.Range("B25").Value = 0.4
$ws.Range("B26").Value = 0
$ws.Range("B27").Value = 0.5
$ws.Range("B28").Value = 0.2132975555746283
$ws.Range("B29").Value = 0.4704848369192122
$ws.Range("B30").Value = 0.3964520361608751
$ws.Range("B31").Value = 0.4230165803625844
$ws.Range("B32").Value = 0.3531481090437554
$ws.Range("B33").Value = 0.4376510431580233
$ws.Range("B34").Value = 0.6485846904589434
$ws.Range("B35").Value = 0.5426591427781329
$ws.Range("B36").Value = 0.2044328674106259
$ws.Range("B37").Value = 0.4167513714050496
$ws.Range("B38").Value = 0.5527780896650785
$ws.Range("B39").Value = 0.4865322308196076
$ws.Range("B40").Value = 0.5234718933384791
$ws.Range("B41").Value = 0.535576229457463
$ws.Range("B42").Value = 0.5969715903930113
$ws.Range("B43").Value = 0.6598534028039785
$ws.Range("B44").Value = 0.6207578930310335
$ws.Range("B45").Value = 0.2938136814264701
$ws.Range("B46").Value = 0
$ws.Range("B47").Value = 0.2964852479966615
$ws.Range("B48").Value = 0.3
$ws.Range("B49").Value = -0.1761138288871404
$ws.Range("B50").Value = -0.02363626827304436
$ws.Range("B51").Value = 0.6810025317521822
$ws.Range("B52").Value = -6.9
$ws.Range("B53").Value = 4.729401638091318
$ws.Range("B54").Value = -0.8905127363963885
$ws.Range("B55").Value = -2.04269378128221
$ws.Range("B56").Value = 2.093024636165651
$ws.Range("B57").Value = 1.098535546956398
$ws.Range("B58").Value = -0.4717175472572421
$ws.Range("B59").Value = 0.9401304606753627
$ws.Range("B60").Value = 0.6392725048137464
$ws.Range("B61").Value = 0.06357296580725347
$ws.Range("B62").Value = -1.012166871044968
$ws.Range("B63").Value = -0.5703626997413522
$ws.Range("B64").Value = 0.2394101325822788
$ws.Range("B65").Value = -0.04072131480353391
$ws.Range("B66").Value = -0.07958838003274593
$ws.Range("B67").Value = 0.02912383308249389
$ws.Range("B68").Value = -0.1311265493919933
$ws.Range("B69").Value = -0.03907468377752821
$ws.Range("B70").Value = 0.1771324545010202
$ws.Range("B71").Value = 0.4946531409412387
$ws.Range("B72").Value = 0.202428137729683
$ws.Range("B73").Value = 0.208573386070384

Write-Host "Edit applied: rows 74-82 deleted, B4:B73 values updated."
